$d = $word.ActiveDocument

function Set-ParaText($index, $newText) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $newText
}

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $old"
    }
}

# ---------------------------------------------------------------------------
# 1. "Registry graph layout Bindings:" block -> "Interfaces (Sets):" block.
#    The old block has 13 list paragraphs; the new block only has 8, so the
#    trailing five (old ISubject/IPredicate/IObject/IStatement/IMapping/
#    ITransform entries minus the one reused as "ITransform : IContext") are
#    deleted outright. Do this first while paragraph numbering is pristine.
# ---------------------------------------------------------------------------
Set-ParaText 68 "Interfaces (Sets):"
Set-ParaText 69 "IContext : ISubject, IPredicate, IObject"
Set-ParaText 70 "ISubject"
Set-ParaText 71 "IPredicate"
Set-ParaText 72 "IObject"
Set-ParaText 73 "IStatement : IContext"
Set-ParaText 74 "IMapping : IContext"
Set-ParaText 75 "ITransform : IContext"

$pStart = $d.Paragraphs.Item(76)
$pEnd = $d.Paragraphs.Item(80)
$rng = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$rng.Delete()

# ---------------------------------------------------------------------------
# 2. Remaining scattered text-only edits (unique substrings, safe to use
#    Find/Replace regardless of the paragraph shift caused by step 1).
# ---------------------------------------------------------------------------
Replace-Text "IObjectKind : IKind, IPredicateResource, ISubjectResource." "IObjectKind : IKind, IPredicateResource, ISubjectResource"
Replace-Text "IStatementKind : IKind, ISubjectKind, IPredicateKind, IObjectKind" "IStatementKind : IKind, IContext, ISubjectKind, IPredicateKind, IObjectKind"
Replace-Text "IMappingKind : IKind, ISubjectKind, IPredicateKind, IObjectKind" "IMappingKind : IKind, IContext, ISubjectKind, IPredicateKind, IObjectKind"
Replace-Text "ITransformKind : IKind, ISubjectKind, IPredicateKind, IObjectKind" "ITransformKind : IKind, IContext, ISubjectKind, IPredicateKind, IObjectKind"

Replace-Text "ISubjectResource : ISubject" "ISubjectResource : IResource, ISubject"
Replace-Text "IPredicateResource : IPredicate" "IPredicateResource : IResource, IPredicate"
Replace-Text "IObjectResource : IObject" "IObjectResource : IResource, IObject"
Replace-Text "IStatementResource : IStatement" "IStatementResource : IResource, IStatement"
Replace-Text "IMappingResource : IMapping" "IMappingResource : IResource, IMapping"
Replace-Text "ITransformResource : ITransform" "ITransformResource : IResource, ITransform"

Replace-Text "ISubjectKind : ISubjectKind" "ISubjectKind"
Replace-Text "IPredicateKind : IPredicateKind" "IPredicateKind"
Replace-Text "(IPredicateOccurrence, IPredicateKind, ISubject, IObject)" "(IPredicateOccurrence, IPredicateKind, ISubjectResource, IObjectResource)"
Replace-Text "IObjectKind : IObjectKind" "IObjectKind"
Replace-Text "(IObjectOccurrence, IObjectKind, IPredicate, ISubject)" "(IObjectOccurrence, IObjectKind, IPredicateResource, ISubjectResource)"
Replace-Text "IStatementKind : IStatementKind" "IStatementKind"
Replace-Text "IMappingKind : IMappingKind" "IMappingKind"
Replace-Text "ITransformKind : ITransformKind" "ITransformKind"

Replace-Text "ISubjectOccurrence : ISubject" "ISubjectOccurrence : IOccurrence, ISubject"
Replace-Text "IPredicateOccurrence : IPredicate" "IPredicateOccurrence : IOccurrence, IPredicate"
Replace-Text "IObjectOccurrence : IObject" "IObjectOccurrence : IOccurrence, IObject"
Replace-Text "IStatementOccurrence : IStatement" "IStatementOccurrence : IOccurrence, IStatement"
Replace-Text "IMappingOccurrence : IMapping" "IMappingOccurrence : IOccurrence, IMapping"
Replace-Text "ITransformOccurrence : ITransform" "ITransformOccurrence : IOccurrence, ITransform"

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
